# Update cryptos list: apply latest price (D) and 1h volume change (E) figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.373.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("D3").Value = "'1.826.80"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'315.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.84%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "'0.4479"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.42%  "
$ws.Range("D8").Value = "'0.3780"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.10%  "
$ws.Range("D9").Value = "'0.07461"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.72%  "
$ws.Range("D10").Value = "'0.8879"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.15%  "
$ws.Range("D11").Value = "'21.04"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("D12").Value = "'1.825.22"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").Value = "'6.757"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.92%  "
$ws.Range("D14").Value = "'5.466"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.14%  "
$ws.Range("D15").Value = "'93.99"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.01%  "
$ws.Range("D16").Value = "'0.07122"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.67%  "
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").Value = "'0.000008808"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").Value = "'15.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.02%  "
$ws.Range("D21").Value = "'27.380.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.96%  "
$ws.Range("D22").Value = "'5.419"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.31%  "
$ws.Range("D23").Value = "'11.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.37%  "
$ws.Range("E24").Value = "  -1.77%  "
$ws.Range("D25").Value = "'151.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("D26").Value = "'2.323"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.67%  "
$ws.Range("D27").Value = "'18.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.21%  "
$ws.Range("D28").Value = "'5.414"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.76%  "
$ws.Range("D29").Value = "'118.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.32%  "
$ws.Range("D30").Value = "'0.08893"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("D31").Value = "'0.7944"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.74%  "
$ws.Range("D32").Value = "'1.207"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.96%  "
$ws.Range("D33").Value = "'4.611"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.14%  "
$ws.Range("D34").Value = "'2.926"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.57%  "
$ws.Range("D35").Value = "'1.000"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  +0.48%  "
$ws.Range("D37").Value = "'0.01992"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.03%  "
$ws.Range("D38").Value = "'0.05316"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("D39").Value = "'7.320"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.69%  "
$ws.Range("D40").Value = "'0.5360"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.82%  "
$ws.Range("D41").Value = "'2.870"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.61%  "
$ws.Range("D42").Value = "'0.1724"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.39%  "
$ws.Range("D43").Value = "'2.308"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +16.01%  "
$ws.Range("D44").Value = "'8.690"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.67%  "
$ws.Range("E45").Value = "  -3.13%  "
$ws.Range("D46").Value = "'10.69"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.47%  "
$ws.Range("D47").Value = "'1.699"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.10%  "
$ws.Range("D48").Value = "'105.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.67%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").Value = "'0.06413"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.66%  "
$ws.Range("D51").Value = "'66.03"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.89%  "
